# Login error validation error message
#
# - Rename "Sheet2" -> "InvalidLogin"
# - Populate InvalidLogin with a Username/Password header row and an
#   "Admin"/"xyz" (invalid password) data row, mirroring ValidLogin's layout
# - Size columns A:B to fit the new content
# - Move the active selection/tab from ValidLogin!B1 to
#   ValidLogin!A1:B2 (no longer the active tab) and make InvalidLogin the
#   active tab with its selection at E3

$wb = $excel.ActiveWorkbook

# --- ValidLogin sheet: sheetView/selection changes only -------------------
$wsValid = $wb.Worksheets.Item("ValidLogin")
[void]$wsValid.Range("A1:B2").Select()

# --- Sheet2 -> InvalidLogin: rename + fill in data -------------------------
$wsInvalid = $wb.Worksheets.Item("Sheet2")
$wsInvalid.Name = "InvalidLogin"

$wsInvalid.Range("A1").Value = "Username"
$wsInvalid.Range("B1").Value = "Password"
$wsInvalid.Range("A2").Value = "Admin"
$wsInvalid.Range("B2").Value = "xyz"

$wsInvalid.Columns.Item(1).ColumnWidth = 10.3333
$wsInvalid.Columns.Item(2).ColumnWidth = 10.3333

# Make InvalidLogin the active sheet/tab and select E3 there
[void]$wsInvalid.Activate()
[void]$wsInvalid.Range("E3").Select()
